$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 84; this shifts the former rows 84-138
# down to rows 85-139, matching the rest of the diff (which is just a
# re-numbering of the existing rows by one).
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 45176
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112012
$ws.Range("G84").Value = "Espinaca"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 40
$ws.Range("K84").Value = 6500
$ws.Range("L84").Value = 6500
$ws.Range("M84").Value = 6500
$ws.Range("N84").Value = "$/cuna 10 kilos"
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 650
$ws.Range("Q84").Value = 10
$ws.Range("R84").Value = "Hortaliza"
